$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.940.37"
$ws.Range("D3").Value = "2.167.81"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'" + "247.24"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "'" + "0.617"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -7.55%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'" + "0.563"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'" + "57.97"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'" + "0.0923"
$ws.Range("E11").Value = "  -5.19%  "
$ws.Range("E12").Value = "  -15.97%  "
$ws.Range("D13").Value = "'" + "0.103"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'" + "14.22"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'" + "0.849"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "2.185.63"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "40.854.44"
$ws.Range("D20").Value = "'" + "0.0" + [char]0x2083 + "0938"
$ws.Range("E20").Value = "  -3.51%  "
$ws.Range("D21").Value = "'" + "71.39"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "'" + "229.13"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D24").Value = "'" + "2.05"
$ws.Range("E24").Value = "  -9.18%  "
$ws.Range("E25").Value = "  +12.51%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "'" + "20.15"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'" + "0.0734"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -9.92%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'" + "59.95"
$ws.Range("E43").Value = "  -14.28%  "
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("E45").Value = "  -11.43%  "
$ws.Range("E46").Value = "  -5.36%  "
$ws.Range("D48").Value = "'" + "0.0987"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("E51").Value = "  -0.91%  "
